$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.427.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.462.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '415.12'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.10'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.765'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +13.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.06'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.47%  '
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000235'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.51%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.75'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.020.22'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.476.92'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.46'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'Polygon'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.09'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '63.270.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '459.76'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.44'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.19%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.37'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.10'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +10.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.29'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.66'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.62%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.38'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.50'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.166'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.112'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '39.90'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.80'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.08'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.82'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.03%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0652'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +59.89%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'LidoDAOToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.33'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.45%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.43'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.318'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.89%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.33'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.83'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.95'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.138'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.95%  '
